$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.412000000000001
$ws.Range("D7").Value = -7.114999999999999
$ws.Range("B8").Value = 6.371
$ws.Range("A12").Value = -21.401
$ws.Range("B12").Value = 6.694999999999999
$ws.Range("B14").Value = 6.532000000000001
$ws.Range("D19").Value = -7.947
$ws.Range("D21").Value = -7.105999999999999
$ws.Range("B22").Value = 6.370999999999999
$ws.Range("D24").Value = -7.456
